$d = $word.ActiveDocument

# Locate the "lib_6_1" text (the DXIL library target-profile string) and
# change the trailing digit from 1 to 3, matching Word's natural
# run-splitting / _GoBack bookmark behavior when a single character is
# replaced in place.
$find = $d.Content
$found = $find.Find.Execute("lib_6_1", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)

if ($found) {
    # The "1" is the last character of the found range.
    $digit = $d.Range($find.End - 1, $find.End)
    $digit.Text = "3"

    # Force a run boundary between "lib_6_" and the new "3" by briefly
    # planting a bookmark right at that seam, then remove it again -
    # Word cannot keep the bookmark inside a single <w:r>, so the split
    # survives even after the temporary bookmark is deleted.
    $seam = $d.Range($find.End - 1, $find.End - 1)
    $d.Bookmarks.Add("zzTempSplit", $seam)

    # Word drops a _GoBack bookmark (collapsed) right after the edited
    # text, marking the last edit location - and since _GoBack is a
    # singleton bookmark, adding it here automatically removes any prior
    # _GoBack bookmark elsewhere in the document.
    $editEnd = $d.Range($find.End, $find.End)
    $d.Bookmarks.Add("_GoBack", $editEnd)

    $d.Bookmarks("zzTempSplit").Delete()
}
